{"js": "// Data dictionary terms consistency edit:\n// 1) Prepend \"Associated with raw data file \" before the filename\n//    \"guadua_leaf_raw_3.xlsx\" in the second paragraph.\n// 2) Change \"large\" -> \"tall\" in the Habit bullet point description.\n\nconst body = context.document.body;\n\n// --- Edit 1: prepend explanatory text before the data file name ---\nconst fileNameResults = body.search(\"guadua_leaf_raw\", { matchCase: true });\nfileNameResults.load(\"items\");\nawait context.sync();\n\nif (fileNameResults.items.length > 0) {\n  fileNameResults.items[0].insertText(\n    \"Associated with raw data file \",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n}\n\n// --- Edit 2: \"large\" -> \"tall\" in the Habit description ---\nconst habitResults = body.search(\"large\", { matchCase: false });\nhabitResults.load(\"items\");\nawait context.sync();\n\nif (habitResults.items.length > 0) {\n  habitResults.items[0].insertText(\"tall\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Data dictionary terms consistency edit:\n# 1) Prepend \"Associated with raw data file \" before the filename\n#    \"guadua_leaf_raw_3.xlsx\" in the second paragraph.\n# 2) Change \"large\" -> \"tall\" in the Habit bullet point description.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: prepend explanatory text before the data file name ---\n$find1 = $d.Content.Find\n$find1.Text = \"guadua_leaf_raw\"\n$find1.Replacement.Text = \"Associated with raw data file guadua_leaf_raw\"\n$find1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- Edit 2: \"large\" -> \"tall\" in the Habit description ---\n$find2 = $d.Content.Find\n$find2.Text = \"large\"\n$find2.Replacement.Text = \"tall\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
